# Append the 2021/09/17 data point (DGS report) as a new row 82
# on the risk-matrix time series sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 82

# Column A holds the report date, stored as text (so it becomes a
# shared-string entry) but displayed with the same date number format
# used by the rest of the column. Temporarily switch the cell to a
# text format while assigning the value so Excel does not silently
# convert the string into a date serial number, then restore the
# date format used by the other cells in the column.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2021/09/17"
$dateCell.NumberFormat = "yyyy/mm/dd"

# Columns B:E hold the numeric series values for the new date.
$ws.Cells.Item($newRow, 2).Value = 173.6
$ws.Cells.Item($newRow, 3).Value = 177.9
$ws.Cells.Item($newRow, 4).Value = 0.83
$ws.Cells.Item($newRow, 5).Value = 0.82

$ws.Range("A" + $newRow).Select()
